{"js": "// Replace each equation cell's text in the single 20x5 table, in row-major\n// order, matching the commit's cell-by-cell numeric updates. Using\n// `cell.value = ...` (TableCell.Text via the Range) preserves the existing\n// run formatting (rFonts/sz) already on each cell, so only the <w:t> content\n// changes -- exactly mirroring the diff.\nconst newValues = [\n  [\n    \"58-24=34\",\n    \"57-11=46\",\n    \"52-3=49\",\n    \"8+7=15\",\n    \"27+47=74\"\n  ],\n  [\n    \"65+23=88\",\n    \"80+6=86\",\n    \"80-47=33\",\n    \"81+9=90\",\n    \"80-23=57\"\n  ],\n  [\n    \"33+51=84\",\n    \"11+23=34\",\n    \"55+12=67\",\n    \"6+13=19\",\n    \"93-23=70\"\n  ],\n  [\n    \"49+44=93\",\n    \"22+14=36\",\n    \"49-23=26\",\n    \"38-22=16\",\n    \"47+52=99\"\n  ],\n  [\n    \"39+50=89\",\n    \"41-37=4\",\n    \"3+82=85\",\n    \"12+7=19\",\n    \"21+50=71\"\n  ],\n  [\n    \"69-43=26\",\n    \"75+19=94\",\n    \"46+14=60\",\n    \"35+31=66\",\n    \"73-57=16\"\n  ],\n  [\n    \"10+62=72\",\n    \"13-4=9\",\n    \"57-10=47\",\n    \"56-26=30\",\n    \"72+1=73\"\n  ],\n  [\n    \"33+28=61\",\n    \"72+5=77\",\n    \"89-43=46\",\n    \"6+83=89\",\n    \"48+33=81\"\n  ],\n  [\n    \"55-41=14\",\n    \"69-12=57\",\n    \"2+70=72\",\n    \"44+12=56\",\n    \"87-65=22\"\n  ],\n  [\n    \"13+50=63\",\n    \"47+5=52\",\n    \"32+38=70\",\n    \"36+55=91\",\n    \"95-43=52\"\n  ],\n  [\n    \"40+54=94\",\n    \"48+8=56\",\n    \"38-28=10\",\n    \"20+66=86\",\n    \"94-37=57\"\n  ],\n  [\n    \"81-69=12\",\n    \"79+14=93\",\n    \"55+19=74\",\n    \"34-0=34\",\n    \"6+61=67\"\n  ],\n  [\n    \"43+8=51\",\n    \"46-29=17\",\n    \"87-25=62\",\n    \"61-18=43\",\n    \"63+19=82\"\n  ],\n  [\n    \"52+19=71\",\n    \"30-28=2\",\n    \"95-39=56\",\n    \"7+11=18\",\n    \"3+18=21\"\n  ],\n  [\n    \"22+49=71\",\n    \"63-28=35\",\n    \"90-13=77\",\n    \"80+0=80\",\n    \"80-16=64\"\n  ],\n  [\n    \"16+33=49\",\n    \"11+12=23\",\n    \"76-40=36\",\n    \"73+20=93\",\n    \"82-8=74\"\n  ],\n  [\n    \"22-5=17\",\n    \"83-13=70\",\n    \"23+9=32\",\n    \"40+57=97\",\n    \"11+81=92\"\n  ],\n  [\n    \"30+50=80\",\n    \"62+18=80\",\n    \"40+50=90\",\n    \"36+55=91\",\n    \"28-9=19\"\n  ],\n  [\n    \"24+46=70\",\n    \"10-0=10\",\n    \"43-39=4\",\n    \"41-22=19\",\n    \"75-0=75\"\n  ],\n  [\n    \"37+55=92\",\n    \"54-52=2\",\n    \"98-97=1\",\n    \"74+10=84\",\n    \"34-9=25\"\n  ]\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(\n    `Expected ${newValues.length} rows, found ${table.rowCount}`\n  );\n}\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(r, c).value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each equation cell's text in the single 20x5 table, in row-major\n# order, matching the commit's cell-by-cell numeric updates. Setting\n# Cell.Range.Text preserves the existing run formatting (rFonts/sz) already\n# on each cell, so only the <w:t> content changes -- exactly mirroring the diff.\n$newValues = @(\n    @(\"58-24=34\", \"57-11=46\", \"52-3=49\", \"8+7=15\", \"27+47=74\"),\n    @(\"65+23=88\", \"80+6=86\", \"80-47=33\", \"81+9=90\", \"80-23=57\"),\n    @(\"33+51=84\", \"11+23=34\", \"55+12=67\", \"6+13=19\", \"93-23=70\"),\n    @(\"49+44=93\", \"22+14=36\", \"49-23=26\", \"38-22=16\", \"47+52=99\"),\n    @(\"39+50=89\", \"41-37=4\", \"3+82=85\", \"12+7=19\", \"21+50=71\"),\n    @(\"69-43=26\", \"75+19=94\", \"46+14=60\", \"35+31=66\", \"73-57=16\"),\n    @(\"10+62=72\", \"13-4=9\", \"57-10=47\", \"56-26=30\", \"72+1=73\"),\n    @(\"33+28=61\", \"72+5=77\", \"89-43=46\", \"6+83=89\", \"48+33=81\"),\n    @(\"55-41=14\", \"69-12=57\", \"2+70=72\", \"44+12=56\", \"87-65=22\"),\n    @(\"13+50=63\", \"47+5=52\", \"32+38=70\", \"36+55=91\", \"95-43=52\"),\n    @(\"40+54=94\", \"48+8=56\", \"38-28=10\", \"20+66=86\", \"94-37=57\"),\n    @(\"81-69=12\", \"79+14=93\", \"55+19=74\", \"34-0=34\", \"6+61=67\"),\n    @(\"43+8=51\", \"46-29=17\", \"87-25=62\", \"61-18=43\", \"63+19=82\"),\n    @(\"52+19=71\", \"30-28=2\", \"95-39=56\", \"7+11=18\", \"3+18=21\"),\n    @(\"22+49=71\", \"63-28=35\", \"90-13=77\", \"80+0=80\", \"80-16=64\"),\n    @(\"16+33=49\", \"11+12=23\", \"76-40=36\", \"73+20=93\", \"82-8=74\"),\n    @(\"22-5=17\", \"83-13=70\", \"23+9=32\", \"40+57=97\", \"11+81=92\"),\n    @(\"30+50=80\", \"62+18=80\", \"40+50=90\", \"36+55=91\", \"28-9=19\"),\n    @(\"24+46=70\", \"10-0=10\", \"43-39=4\", \"41-22=19\", \"75-0=75\"),\n    @(\"37+55=92\", \"54-52=2\", \"98-97=1\", \"74+10=84\", \"34-9=25\"),\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nif ($tbl.Rows.Count -ne $newValues.Count) {\n    throw \"Expected $($newValues.Count) rows, found $($tbl.Rows.Count)\"\n}\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    for ($c = 0; $c -lt $newValues[$r].Count; $c++) {\n        $cell = $tbl.Cell($r + 1, $c + 1)\n        $cell.Range.Text = $newValues[$r][$c]\n    }\n}\n"}
